$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New HIGH/LOW/CLOSE/LTP/VOL/9:25-CLOSE figures per row (row => col => value)
$updates = @{
    2 = @{ "B" = 2465.55; "C" = 2395.6; "D" = 2449.9; "E" = 2450.05; "F" = 67; "G" = 2401.95 }
    3 = @{ "B" = 390.55; "C" = 383.8; "D" = 385.15; "E" = 384.75; "F" = 26; "G" = 387.95 }
    4 = @{ "B" = 1518.7; "C" = 1490.15; "D" = 1514; "E" = 1514.5; "F" = 9; "G" = 1493.25 }
    5 = @{ "B" = 7329.4; "C" = 7155; "D" = 7300.1; "E" = 7320.6; "F" = 6; "G" = 7185 }
    6 = @{ "B" = 236.25; "C" = 231.6; "D" = 235.15; "E" = 235.1; "F" = 81; "G" = 232.4 }
    7 = @{ "B" = 191.8; "C" = 186.3; "D" = 190.75; "E" = 190.6; "F" = 158; "G" = 187.6 }
    8 = @{ "B" = 237.25; "C" = 232.1; "D" = 236.9; "E" = 236.9; "F" = 148; "G" = 232.4 }
    9 = @{ "B" = 509.5; "C" = 500.3; "D" = 503.9; "E" = 503.25; "F" = 30; "G" = 504.15 }
    10 = @{ "B" = 3406.45; "C" = 3322.1; "D" = 3402; "E" = 3400.4; "F" = 6; "G" = 3330.05 }
    11 = @{ "B" = 144.65; "C" = 141.3; "D" = 143.5; "E" = 143.8; "F" = 136; "G" = 143.7 }
    12 = @{ "B" = 1188; "C" = 1171.1; "D" = 1185.1; "E" = 1185.8; "F" = 16; "G" = 1172.1 }
    13 = @{ "B" = 1579; "C" = 1559.35; "D" = 1573.75; "E" = 1574.7; "F" = 216; "G" = 1571.35 }
    14 = @{ "B" = 478.5; "C" = 468.75; "D" = 472.7; "E" = 473.25; "F" = 111; "G" = 469.75 }
    15 = @{ "B" = 973.05; "C" = 952.8; "D" = 967.6; "E" = 968.7; "F" = 162; "G" = 955.6 }
    16 = @{ "B" = 1423.6; "C" = 1376.3; "D" = 1415; "E" = 1416.55; "F" = 27; "G" = 1383.35 }
    17 = @{ "B" = 1446; "C" = 1428.2; "D" = 1442; "E" = 1443.85; "F" = 35; "G" = 1429.55 }
    18 = @{ "B" = 713.7; "C" = 691.05; "D" = 701; "E" = 701.1; "F" = 38; "G" = 692.95 }
    19 = @{ "B" = 435.85; "C" = 425.2; "D" = 434.1; "E" = 433.2; "F" = 65; "G" = 425.7 }
    20 = @{ "B" = 1604.5; "C" = 1582.4; "D" = 1589.95; "E" = 1591.95; "F" = 21; "G" = 1599.05 }
    21 = @{ "B" = 300.2; "C" = 293.85; "D" = 298.6; "E" = 299.15; "F" = 22; "G" = 296.6 }
    22 = @{ "B" = 2425.65; "C" = 2402.1; "D" = 2411.4; "E" = 2412.65; "F" = 90; "G" = 2411.85 }
    23 = @{ "B" = 571.2; "C" = 562.3; "D" = 569.5; "E" = 569.55; "F" = 181; "G" = 564.5 }
    24 = @{ "B" = 623.75; "C" = 615.35; "D" = 619.55; "E" = 620.8; "F" = 9; "G" = 617.8 }
    25 = @{ "B" = 1076; "C" = 1064; "D" = 1071; "E" = 1071.15; "F" = 5; "G" = 1068.25 }
    26 = @{ "B" = 614.9; "C" = 603.7; "D" = 610.4; "E" = 611.2; "F" = 102; "G" = 604.05 }
    27 = @{ "B" = 255.95; "C" = 247.6; "D" = 255.5; "E" = 255.35; "F" = 280; "G" = 247.75 }
    28 = @{ "B" = 128.7; "C" = 125.4; "D" = 126.8; "E" = 127.05; "F" = 764; "G" = 125.7 }
    29 = @{ "B" = 8339.799999999999; "C" = 8249; "D" = 8269.9; "E" = 8262.4; "G" = 8312.1 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
